$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(101).Insert()

$ws.Range("A101").Value = 3
$ws.Range("B101").Value = "Femacal de La Calera"
$ws.Range("C101").Value = "Coquimbo"
$ws.Range("D101").Value = 44662
$ws.Range("E101").Value = 5
$ws.Range("F101").Value = 100112039
$ws.Range("G101").Value = "Ciboulette"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 160
$ws.Range("K101").Value = 1500
$ws.Range("L101").Value = 1500
$ws.Range("M101").Value = 1500
$ws.Range("N101").Value = "$/docena de atados"
$ws.Range("O101").Value = "Provincia de Quillota"
$ws.Range("P101").Value = 500
$ws.Range("Q101").Value = 3
$ws.Range("R101").Value = "Hortaliza"
